# AutoSPInstaller_Issue_Tracker.xlsx - update for Version 3.86
# - Adds trailing periods to several existing Issue/Work around/Solution cells
# - Normalizes D2:F3 text-cell formatting (vertical-top + wrap text, drop explicit
#   left-horizontal alignment so the duplicate style gets collapsed)
# - Appends a new issue row (row 4) for SharePoint 2013 / AutoSPInstaller 3.86 / ChangeSet 99664

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append trailing periods to existing text that is missing them ---

$ws.Range("D2").Value = "When you specify a 'non default' search index location path in your 'IndexLocation' parameter in the 'AutoSPinstallerInput' XML file you get the following crawl log error: ""The filtering process has been terminated""."
$ws.Range("E2").Value = "In the 'AutoSPInstallerInput'  XML file under the 'EnterpriseSearchService' section leave the default location for the search index in place: IndexLocation=""C:\Program Files\Microsoft Office Servers\15.0\Data\Office Server\Applications""."
$ws.Range("F2").Value = "Meant to have been resolved in changeset 99077."

$ws.Range("D3").Value = "When you launch the User Profile Service Application and attempt to edit your 'Synchronization Connections' created during the install; you get a 'Unable to process Put message' exception when attempting to save the changes."
$ws.Range("E3").Value = "Creating a whole new Synchronization Connection in the User Profile Service Application under 'Synchronization' --> 'Configure Synchronization Connections' should resolve this. Check and test this with the Forefront Synchronization Service Manager (FIM 2010) client that gets installed with SharePoint. This is the default location for the FIM client: C:\Program Files\Microsoft Office Servers\15.0\Synchronization Service\UIShell\miisclient.exe."

# --- 2. Normalize formatting on the text columns (D:F) for the existing data rows ---
# vertical-top + wrap text, no explicit horizontal alignment.
# NB: only touch the cells that actually hold content (D2,E2,F2,D3,E3) - row 3 has
# no F3 cell, so a contiguous D2:F3 range would create a phantom blank F3 cell.
# NB: multi-area ranges ("A1,B2") only apply formatting to the first area here, so
# loop over each address individually instead.

foreach ($addr in @("D2","E2","F2","D3","E3")) {
    $r = $ws.Range($addr)
    $r.WrapText = $true
    $r.VerticalAlignment = -4160   # xlVAlignTop
    $r.HorizontalAlignment = 1     # xlGeneral
}

# --- 3. Add the new issue row (row 4) ---

$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = 3.86
$ws.Range("C4").Value = 99664
$ws.Range("D4").Value = "When you run the 'AutoSPInstallerLaunch' BAT file during the 'PrerequisiteInstallerFiles' phase you get .Net 3.5.1 framework installation exception messages. Even when you change the '<OfflineInstall>true</OfflineInstall>' value to 'false' in in the 'AutoSPinstallerInput' XML file; you still encounter issues with installing the framework."
$ws.Range("F4").Value = "Copy all of the files from your Windows Server 2012 Installation media '\sources\sxs' location to your AutoSPInstaller  sxs directory '\SP\2013\SharePoint\PrerequisiteInstallerFiles\sxs'. Try running the AutoSPInstallerLaunch' BAT file again."

# Match the vertical-top alignment used by the other numeric cells in A:C
$ws.Range("A4:C4").VerticalAlignment = -4160   # xlVAlignTop

# Row 4 has no E4 cell, so format D4 and F4 individually (not a contiguous D4:F4
# range) to avoid materializing a phantom blank E4 cell.
foreach ($addr in @("D4","F4")) {
    $r = $ws.Range($addr)
    $r.WrapText = $true
    $r.VerticalAlignment = -4160   # xlVAlignTop
    $r.HorizontalAlignment = 1     # xlGeneral
}

$ws.Rows.Item(4).RowHeight = 180
